$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.179.10"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.828.91"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.09"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6157"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07339"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2905"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.16"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07635"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.832.71"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.977"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6702"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.43"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008967"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.841"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "29.169.56"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "2.078.68"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.32"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.364"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.67"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.512"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  +4.62%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.087"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.847"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7190"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.614"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.860"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("D39").Value = "1.223.33"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01759"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.176"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8995"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "1.998.26"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.90"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.35"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5039"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.188"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000117"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4023"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  +5.20%  "
